$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.609230666666667
$ws.Cells.Item(2, 8).Value = 4.827692
$ws.Cells.Item(2, 9).Value = 0.5482851650894511
$ws.Cells.Item(2, 10).Value = 0.5482851650894512
$ws.Cells.Item(2, 13).Value = 1.174933333333333
$ws.Cells.Item(2, 14).Value = 3.5248
$ws.Cells.Item(2, 15).Value = 0.01171850713626266
$ws.Cells.Item(2, 16).Value = 0.01171850713626266
$ws.Cells.Item(2, 17).Value = 1.890738751288889
$ws.Cells.Item(2, 18).Value = 17.0166487616
$ws.Cells.Item(2, 19).Value = 0.006425083619807682
$ws.Cells.Item(2, 20).Value = 0.006425083619807682
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.609230666666667
$ws.Cells.Item(3, 8).Value = 4.827692
$ws.Cells.Item(3, 9).Value = 0.5482851650894511
$ws.Cells.Item(3, 10).Value = 0.5482851650894512
$ws.Cells.Item(3, 15).Value = 0.2743256641287217
$ws.Cells.Item(3, 16).Value = 0.2743256641287218
$ws.Cells.Item(3, 17).Value = 44.26145392156622
$ws.Cells.Item(3, 18).Value = 398.353085294096
$ws.Cells.Item(3, 19).Value = 0.1504086920450895
$ws.Cells.Item(3, 20).Value = 0.1504086920450896
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.609230666666667
$ws.Cells.Item(4, 8).Value = 4.827692
$ws.Cells.Item(4, 9).Value = 0.5482851650894511
$ws.Cells.Item(4, 10).Value = 0.5482851650894512
$ws.Cells.Item(4, 13).Value = 39.361408
$ws.Cells.Item(4, 14).Value = 118.084224
$ws.Cells.Item(4, 15).Value = 0.3925813724534833
$ws.Cells.Item(4, 16).Value = 0.3925813724534833
$ws.Cells.Item(4, 17).Value = 63.34158483677867
$ws.Cells.Item(4, 18).Value = 570.074263531008
$ws.Cells.Item(4, 19).Value = 0.2152465426067014
$ws.Cells.Item(4, 20).Value = 0.2152465426067015
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.609230666666667
$ws.Cells.Item(5, 8).Value = 4.827692
$ws.Cells.Item(5, 9).Value = 0.5482851650894511
$ws.Cells.Item(5, 10).Value = 0.5482851650894512
$ws.Cells.Item(5, 13).Value = 32.221985
$ws.Cells.Item(5, 14).Value = 96.665955
$ws.Cells.Item(5, 15).Value = 0.3213744562815322
$ws.Cells.Item(5, 16).Value = 0.3213744562815322
$ws.Cells.Item(5, 17).Value = 51.85260640287333
$ws.Cells.Item(5, 18).Value = 466.67345762586
$ws.Cells.Item(5, 19).Value = 0.1762048468178525
$ws.Cells.Item(5, 20).Value = 0.1762048468178525
$ws.Cells.Item(6, 7).Value = 0.9591736666666666
$ws.Cells.Item(6, 9).Value = 0.3268025542087943
$ws.Cells.Item(6, 10).Value = 0.3268025542087943
$ws.Cells.Item(6, 13).Value = 1.174933333333333
$ws.Cells.Item(6, 14).Value = 3.5248
$ws.Cells.Item(6, 15).Value = 0.01171850713626266
$ws.Cells.Item(6, 16).Value = 0.01171850713626266
$ws.Cells.Item(6, 17).Value = 1.126965113422222
$ws.Cells.Item(6, 18).Value = 10.1426860208
$ws.Cells.Item(6, 19).Value = 0.00382963806364462
$ws.Cells.Item(6, 20).Value = 0.00382963806364462
$ws.Cells.Item(7, 7).Value = 0.9591736666666666
$ws.Cells.Item(7, 9).Value = 0.3268025542087943
$ws.Cells.Item(7, 10).Value = 0.3268025542087943
$ws.Cells.Item(7, 15).Value = 0.2743256641287217
$ws.Cells.Item(7, 16).Value = 0.2743256641287218
$ws.Cells.Item(7, 17).Value = 26.38181208532755
$ws.Cells.Item(7, 19).Value = 0.08965032772229008
$ws.Cells.Item(7, 20).Value = 0.0896503277222901
$ws.Cells.Item(8, 7).Value = 0.9591736666666666
$ws.Cells.Item(8, 9).Value = 0.3268025542087943
$ws.Cells.Item(8, 10).Value = 0.3268025542087943
$ws.Cells.Item(8, 13).Value = 39.361408
$ws.Cells.Item(8, 14).Value = 118.084224
$ws.Cells.Item(8, 15).Value = 0.3925813724534833
$ws.Cells.Item(8, 16).Value = 0.3925813724534833
$ws.Cells.Item(8, 17).Value = 37.75442603652267
$ws.Cells.Item(8, 18).Value = 339.789834328704
$ws.Cells.Item(8, 19).Value = 0.1282965952525923
$ws.Cells.Item(8, 20).Value = 0.1282965952525924
$ws.Cells.Item(9, 7).Value = 0.9591736666666666
$ws.Cells.Item(9, 9).Value = 0.3268025542087943
$ws.Cells.Item(9, 10).Value = 0.3268025542087943
$ws.Cells.Item(9, 13).Value = 32.221985
$ws.Cells.Item(9, 14).Value = 96.665955
$ws.Cells.Item(9, 15).Value = 0.3213744562815322
$ws.Cells.Item(9, 16).Value = 0.3213744562815322
$ws.Cells.Item(9, 17).Value = 30.90647949972833
$ws.Cells.Item(9, 18).Value = 278.158315497555
$ws.Cells.Item(9, 19).Value = 0.1050259931702672
$ws.Cells.Item(9, 20).Value = 0.1050259931702672
$ws.Cells.Item(10, 9).Value = 0.04674417878325851
$ws.Cells.Item(10, 10).Value = 0.04674417878325852
$ws.Cells.Item(10, 13).Value = 1.174933333333333
$ws.Cells.Item(10, 14).Value = 3.5248
$ws.Cells.Item(10, 15).Value = 0.01171850713626266
$ws.Cells.Item(10, 16).Value = 0.01171850713626266
$ws.Cells.Item(10, 17).Value = 0.1611953703111111
$ws.Cells.Item(10, 18).Value = 1.4507583328
$ws.Cells.Item(10, 19).Value = 0.0005477719926503524
$ws.Cells.Item(10, 20).Value = 0.0005477719926503524
$ws.Cells.Item(11, 9).Value = 0.04674417878325851
$ws.Cells.Item(11, 10).Value = 0.04674417878325852
$ws.Cells.Item(11, 15).Value = 0.2743256641287217
$ws.Cells.Item(11, 16).Value = 0.2743256641287218
$ws.Cells.Item(11, 17).Value = 3.773520509129778
$ws.Cells.Item(11, 19).Value = 0.01282312788886909
$ws.Cells.Item(11, 20).Value = 0.0128231278888691
$ws.Cells.Item(12, 9).Value = 0.04674417878325851
$ws.Cells.Item(12, 10).Value = 0.04674417878325852
$ws.Cells.Item(12, 13).Value = 39.361408
$ws.Cells.Item(12, 14).Value = 118.084224
$ws.Cells.Item(12, 15).Value = 0.3925813724534833
$ws.Cells.Item(12, 16).Value = 0.3925813724534833
$ws.Cells.Item(12, 17).Value = 5.400201491029334
$ws.Cells.Item(12, 18).Value = 48.601813419264
$ws.Cells.Item(12, 19).Value = 0.01835089386094262
$ws.Cells.Item(12, 20).Value = 0.01835089386094263
$ws.Cells.Item(13, 9).Value = 0.04674417878325851
$ws.Cells.Item(13, 10).Value = 0.04674417878325852
$ws.Cells.Item(13, 13).Value = 32.221985
$ws.Cells.Item(13, 14).Value = 96.665955
$ws.Cells.Item(13, 15).Value = 0.3213744562815322
$ws.Cells.Item(13, 16).Value = 0.3213744562815322
$ws.Cells.Item(13, 17).Value = 4.420705972736666
$ws.Cells.Item(13, 18).Value = 39.78635375463
$ws.Cells.Item(13, 19).Value = 0.01502238504079643
$ws.Cells.Item(13, 20).Value = 0.01502238504079644
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.2294253333333333
$ws.Cells.Item(14, 8).Value = 0.688276
$ws.Cells.Item(14, 9).Value = 0.07816810191849585
$ws.Cells.Item(14, 10).Value = 0.07816810191849587
$ws.Cells.Item(14, 13).Value = 1.174933333333333
$ws.Cells.Item(14, 14).Value = 3.5248
$ws.Cells.Item(14, 15).Value = 0.01171850713626266
$ws.Cells.Item(14, 16).Value = 0.01171850713626266
$ws.Cells.Item(14, 17).Value = 0.2695594716444445
$ws.Cells.Item(14, 18).Value = 2.4260352448
$ws.Cells.Item(14, 19).Value = 0.0009160134601600003
$ws.Cells.Item(14, 20).Value = 0.0009160134601600004
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.2294253333333333
$ws.Cells.Item(15, 8).Value = 0.688276
$ws.Cells.Item(15, 9).Value = 0.07816810191849585
$ws.Cells.Item(15, 10).Value = 0.07816810191849587
$ws.Cells.Item(15, 15).Value = 0.2743256641287217
$ws.Cells.Item(15, 16).Value = 0.2743256641287218
$ws.Cells.Item(15, 17).Value = 6.310281695543112
$ws.Cells.Item(15, 18).Value = 56.79253525988801
$ws.Cells.Item(15, 19).Value = 0.02144351647247298
$ws.Cells.Item(15, 20).Value = 0.02144351647247299
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.2294253333333333
$ws.Cells.Item(16, 8).Value = 0.688276
$ws.Cells.Item(16, 9).Value = 0.07816810191849585
$ws.Cells.Item(16, 10).Value = 0.07816810191849587
$ws.Cells.Item(16, 13).Value = 39.361408
$ws.Cells.Item(16, 14).Value = 118.084224
$ws.Cells.Item(16, 15).Value = 0.3925813724534833
$ws.Cells.Item(16, 16).Value = 0.3925813724534833
$ws.Cells.Item(16, 17).Value = 9.030504150869335
$ws.Cells.Item(16, 18).Value = 81.274537357824
$ws.Cells.Item(16, 19).Value = 0.03068734073324687
$ws.Cells.Item(16, 20).Value = 0.03068734073324687
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.2294253333333333
$ws.Cells.Item(17, 8).Value = 0.688276
$ws.Cells.Item(17, 9).Value = 0.07816810191849585
$ws.Cells.Item(17, 10).Value = 0.07816810191849587
$ws.Cells.Item(17, 13).Value = 32.221985
$ws.Cells.Item(17, 14).Value = 96.665955
$ws.Cells.Item(17, 15).Value = 0.3213744562815322
$ws.Cells.Item(17, 16).Value = 0.3213744562815322
$ws.Cells.Item(17, 17).Value = 7.392539649286666
$ws.Cells.Item(17, 18).Value = 66.53285684357999
$ws.Cells.Item(17, 19).Value = 0.02512123125261599
$ws.Cells.Item(17, 20).Value = 0.02512123125261601
